$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.363110780715942
$ws.Range("B1").Value = 2.153929233551025
$ws.Range("C1").Value = 4.866361141204834
$ws.Range("D1").Value = 3.454369068145752
$ws.Range("E1").Value = 1.266086101531982
